$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.238.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.78%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.606.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.48%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("E5").Value = "  +0.01%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3763"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "53.12"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.28%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3612"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.62%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.265"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08145"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.02%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.88"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.608"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.67%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.352"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.24%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001243"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.22%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.602.55"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.91%  "

# Row 19
$ws.Range("E19").Value = "  +1.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.16"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.10%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.524"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.39%  "

# Row 22
$ws.Range("E22").Value = "  +0.19%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.225.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.70%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.086"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.21%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.408"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.45%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.44%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "150.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.259"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.01%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.32%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.410"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.759"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.780.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.27%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9494"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02767"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.14%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.30"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.18%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.07388"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.38%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2513"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.114"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.31%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08748"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.81%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.397"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.98%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7097"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.94%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.80%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.92%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6535"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.326"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.53%  "

# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.04%  "

# Row 48
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.25%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "134.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.48%  "

# Row 50
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07968"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.54%  "

# Row 51
$ws.Range("B51").Value = "Flow"
$ws.Range("C51").Value = "https://coinranking.com/coin/QQ0NCmjVq+flow-flow"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.197"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
